$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 7
$ws.Cells.Item($row, 1).Value = "Mateusz"
$ws.Cells.Item($row, 2).Value = "Lugowski"

# Force the phone number to be stored as text (shared string) rather than
# being auto-converted to a number, mirroring the other phone/id columns
# in the sheet (e.g. "03222222222") that keep their literal digits. Clear
# the transient formatting afterwards so the cell keeps the workbook's
# default (unstyled) appearance.
$c = $ws.Cells.Item($row, 3)
$c.NumberFormat = "@"
$c.Value = "11111111111"
$c.ClearFormats()

$ws.Cells.Item($row, 4).Value = "lugowski.mateusz.02@gmail.com"
